# Auto-generated Excel COM-interop script to update F-column popularity values
# across the four worksheets of the workbook, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 225
$ws.Cells.Item(4, 6).Value = 460
$ws.Cells.Item(5, 6).Value = 1934
$ws.Cells.Item(7, 6).Value = 7630
$ws.Cells.Item(8, 6).Value = 226
$ws.Cells.Item(9, 6).Value = 6
$ws.Cells.Item(11, 6).Value = 213
$ws.Cells.Item(12, 6).Value = 1719
$ws.Cells.Item(13, 6).Value = 1478
$ws.Cells.Item(14, 6).Value = 1290
$ws.Cells.Item(16, 6).Value = 3657
$ws.Cells.Item(17, 6).Value = 5888
$ws.Cells.Item(18, 6).Value = 644
$ws.Cells.Item(19, 6).Value = 6
$ws.Cells.Item(20, 6).Value = 1021
$ws.Cells.Item(21, 6).Value = 1208
$ws.Cells.Item(22, 6).Value = 384
$ws.Cells.Item(23, 6).Value = 5989
$ws.Cells.Item(26, 6).Value = 4074
$ws.Cells.Item(27, 6).Value = 233
$ws.Cells.Item(28, 6).Value = 678
$ws.Cells.Item(29, 6).Value = 1879
$ws.Cells.Item(30, 6).Value = 1134
$ws.Cells.Item(31, 6).Value = 271
$ws.Cells.Item(32, 6).Value = 3
$ws.Cells.Item(33, 6).Value = 22
$ws.Cells.Item(34, 6).Value = 176
$ws.Cells.Item(35, 6).Value = 312
$ws.Cells.Item(36, 6).Value = 1124
$ws.Cells.Item(37, 6).Value = 483
$ws.Cells.Item(38, 6).Value = 1836
$ws.Cells.Item(40, 6).Value = 376
$ws.Cells.Item(41, 6).Value = 140
$ws.Cells.Item(42, 6).Value = 1069
$ws.Cells.Item(47, 6).Value = 74
$ws.Cells.Item(48, 6).Value = 154
$ws.Cells.Item(49, 6).Value = 14

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 118
$ws.Cells.Item(10, 6).Value = 653
$ws.Cells.Item(11, 6).Value = 338
$ws.Cells.Item(14, 6).Value = 196
$ws.Cells.Item(15, 6).Value = 102
$ws.Cells.Item(18, 6).Value = 340
$ws.Cells.Item(19, 6).Value = 145
$ws.Cells.Item(25, 6).Value = 122

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 3325
$ws.Cells.Item(5, 6).Value = 441
$ws.Cells.Item(7, 6).Value = 1535
$ws.Cells.Item(8, 6).Value = 460
$ws.Cells.Item(9, 6).Value = 3042
$ws.Cells.Item(10, 6).Value = 389
$ws.Cells.Item(11, 6).Value = 831
$ws.Cells.Item(12, 6).Value = 984
$ws.Cells.Item(13, 6).Value = 1097
$ws.Cells.Item(14, 6).Value = 1487

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1535
$ws.Cells.Item(3, 6).Value = 460
$ws.Cells.Item(4, 6).Value = 460
$ws.Cells.Item(5, 6).Value = 3042
$ws.Cells.Item(6, 6).Value = 1934
$ws.Cells.Item(8, 6).Value = 7630
$ws.Cells.Item(9, 6).Value = 226
$ws.Cells.Item(12, 6).Value = 1719
$ws.Cells.Item(13, 6).Value = 1478
$ws.Cells.Item(14, 6).Value = 1097
$ws.Cells.Item(15, 6).Value = 1290
$ws.Cells.Item(16, 6).Value = 653
$ws.Cells.Item(18, 6).Value = 1487
$ws.Cells.Item(19, 6).Value = 3657
$ws.Cells.Item(20, 6).Value = 338
$ws.Cells.Item(22, 6).Value = 644
$ws.Cells.Item(23, 6).Value = 1021
$ws.Cells.Item(24, 6).Value = 1208
$ws.Cells.Item(25, 6).Value = 384
$ws.Cells.Item(26, 6).Value = 5990
$ws.Cells.Item(28, 6).Value = 678
$ws.Cells.Item(29, 6).Value = 1879
$ws.Cells.Item(30, 6).Value = 1134
$ws.Cells.Item(31, 6).Value = 271
$ws.Cells.Item(32, 6).Value = 22
$ws.Cells.Item(33, 6).Value = 145
$ws.Cells.Item(34, 6).Value = 176
$ws.Cells.Item(35, 6).Value = 312
$ws.Cells.Item(36, 6).Value = 1124
$ws.Cells.Item(37, 6).Value = 483
$ws.Cells.Item(38, 6).Value = 1836
$ws.Cells.Item(41, 6).Value = 140
$ws.Cells.Item(42, 6).Value = 1069
$ws.Cells.Item(43, 6).Value = 122
$ws.Cells.Item(46, 6).Value = 74
$ws.Cells.Item(48, 6).Value = 154
$ws.Cells.Item(49, 6).Value = 14
